# Updated symbol list on Mon Jan  9 20:39:36 UTC 2023 with GitHub Actions
#
# Applies the refreshed crypto price / volume(1h) figures (and the
# FTXToken <-> GateToken row swap) to the "cryptos" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (avoids Excel auto-converting
# numeric-looking / percent-looking strings into real numbers), while
# not leaving behind any extra number-format / style residue.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2 - BNB
Set-TextValue "D2" "274.48"
Set-TextValue "E2" "2.09%"

# Row 3 - OKB
Set-TextValue "D3" "26.80"
Set-TextValue "E3" "0.38%"

# Row 4 - HuobiToken
Set-TextValue "D4" "4.917"
Set-TextValue "E4" "4.40%"

# Row 5 - Cronos
Set-TextValue "D5" "0.06345"
Set-TextValue "E5" "4.01%"

# Row 6 - KuCoinToken
Set-TextValue "D6" "6.951"
Set-TextValue "E6" "3.17%"

# Row 7 - was FTXToken, now GateToken
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "3.357"
Set-TextValue "E7" "6.01%"

# Row 8 - was GateToken, now FTXToken
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D8" "1.480"
Set-TextValue "E8" "65.81%"

# Row 9 - MXToken
Set-TextValue "D9" "0.8879"
Set-TextValue "E9" "3.49%"

# Row 10 - WazirX
Set-TextValue "D10" "0.1472"
Set-TextValue "E10" "3.88%"

# Row 11 - LiechtensteinCryptoassetsExchange
Set-TextValue "D11" "0.05166"
Set-TextValue "E11" "3.41%"

# Row 12 - MandalaExchangeToken
Set-TextValue "D12" "0.07413"
Set-TextValue "E12" "4.48%"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.03171"
Set-TextValue "E13" "-1.39%"

# Row 14 - BitMartToken
Set-TextValue "D14" "0.09057"
Set-TextValue "E14" "0.28%"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001576"
Set-TextValue "E15" "2.89%"

# Row 16 - One
Set-TextValue "D16" "0.0006310"
Set-TextValue "E16" "3.92%"

# Row 17 - TigerCash
Set-TextValue "D17" "0.006036"
Set-TextValue "E17" "-0.24%"

# Row 18 - LEO
Set-TextValue "D18" "3.483"
Set-TextValue "E18" "0.63%"

# Row 19 - BTSEToken (price unchanged, only volume)
Set-TextValue "E19" "1.65%"

# Row 21 - ProBitToken
Set-TextValue "D21" "0.1335"
Set-TextValue "E21" "2.71%"

# Row 22 - MCDex
Set-TextValue "D22" "3.940"
Set-TextValue "E22" "2.64%"

# Row 23 - CoinExToken
Set-TextValue "D23" "0.04330"
Set-TextValue "E23" "2.32%"

# Row 24 - BitKan
Set-TextValue "D24" "0.001182"
Set-TextValue "E24" "-0.36%"

# Row 25 - HotbitToken
Set-TextValue "D25" "0.003658"
Set-TextValue "E25" "-11.79%"

# Row 26 - NitroEx (price unchanged, only volume)
Set-TextValue "E26" "0.27%"

# Row 27
Set-TextValue "D27" "0.0001943"
Set-TextValue "E27" "15.55%"

# Row 40 - IDEX
Set-TextValue "D40" "0.04033"
Set-TextValue "E40" "2.13%"

# Row 41 - KickToken
Set-TextValue "D41" "0.006629"
Set-TextValue "E41" "58.51%"

# Row 42 - BKEXToken
Set-TextValue "D42" "0.1167"
Set-TextValue "E42" "4.35%"

# Row 43 - CEJI
Set-TextValue "D43" "0.002368"
Set-TextValue "E43" "20.76%"

# Row 44 - LocalTraders (price unchanged, only volume)
Set-TextValue "E44" "-2.67%"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005250"
Set-TextValue "E45" "2.11%"

# Row 46 - BOLO (price unchanged, only volume)
Set-TextValue "E46" "389.66%"

# Row 47 - CoinbaseStockToken
Set-TextValue "D47" "0.02127"
Set-TextValue "E47" "-13.10%"

# Row 48 - SpecialPowerGold (price unchanged, only volume)
Set-TextValue "E48" "-0.12%"
